# "method for regular proposal started"
# Rename the special-proposal labels to the new "SPS" naming, and move the
# active selection on Sheet1 from B2 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the proposal name cells (A2/B2 on Sheet1).
$ws.Range("A2").Value = "NEW-SPS"
$ws.Range("B2").Value = "NEW SPS"

# Move the active selection to C3 (was B2).
$ws.Activate()
$ws.Range("C3").Select()
